$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.438467860221863
$ws.Range("B1").Value = 3.517448663711548
$ws.Range("C1").Value = 5.267248630523682
$ws.Range("D1").Value = 1.723904609680176
$ws.Range("E1").Value = 0.9659792184829712
